$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22..118 down to 23..119
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly price entry.
# (Columns A,B,C,E,F,G,H,I,J,K,L,Q,T carry the same constant metadata as the
# other rows for this market/product; D,M,N,O,P,R,S are the row-specific values.)
$ws.Cells.Item(22, 1).Value = 11
$ws.Cells.Item(22, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(22, 3).Value = "Bíobío"
$ws.Cells.Item(22, 4).Value = 44764
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100108
$ws.Cells.Item(22, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(22, 9).Value = 100108002
$ws.Cells.Item(22, 10).Value = "Mango"
$ws.Cells.Item(22, 11).Value = "Sin especificar"
$ws.Cells.Item(22, 12).Value = "Primera"
$ws.Cells.Item(22, 13).Value = 200
$ws.Cells.Item(22, 14).Value = 8000
$ws.Cells.Item(22, 15).Value = 8500
$ws.Cells.Item(22, 16).Value = 8250
$ws.Cells.Item(22, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(22, 18).Value = "Brasil"
$ws.Cells.Item(22, 19).Value = 2062
$ws.Cells.Item(22, 20).Value = 4
